# Fix incorrect IFRS financial figures for company_list rows 2-9 (data entry error)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 21841
$ws.Range("E2").Value = 930
$ws.Range("F2").Value = 1023
$ws.Range("G2").Value = 325
$ws.Range("H2").Value = 206
$ws.Range("I2").Value = 187
$ws.Range("J2").Value = 19
$ws.Range("K2").Value = 41000
$ws.Range("L2").Value = 17834
$ws.Range("M2").Value = 23166
$ws.Range("N2").Value = 22737
$ws.Range("O2").Value = 429
$ws.Range("P2").Value = 68
$ws.Range("Q2").Value = 1450
$ws.Range("R2").Value = -3867
$ws.Range("S2").Value = 2309
$ws.Range("T2").Value = 2944
$ws.Range("U2").Value = -1494
$ws.Range("V2").Value = 9780
$ws.Range("W2").Value = 4.26
$ws.Range("X2").Value = 0.9399999999999999
$ws.Range("Y2").Value = 0.8
$ws.Range("Z2").Value = 0.51
$ws.Range("AA2").Value = 76.98
$ws.Range("AB2").Value = 27977.11
$ws.Range("AC2").Value = 1376
$ws.Range("AE2").Value = 167528
$ws.Range("AF2").Value = 8.779999999999999
$ws.Range("AG2").Value = 490
$ws.Range("AH2").Value = 0.03
$ws.Range("AI2").Value = 35.63
$ws.Range("AJ2").Value = 12372030
$ws.Range("AD2").ClearContents()  # column removed for this row after the fix

# Row 3
$ws.Range("D3").Value = 22992
$ws.Range("E3").Value = 1429
$ws.Range("F3").Value = 1429
$ws.Range("G3").Value = 979
$ws.Range("H3").Value = 1003
$ws.Range("I3").Value = 996
$ws.Range("J3").Value = 7
$ws.Range("K3").Value = 43820
$ws.Range("L3").Value = 20194
$ws.Range("M3").Value = 23626
$ws.Range("N3").Value = 23361
$ws.Range("O3").Value = 265
$ws.Range("P3").Value = 68
$ws.Range("Q3").Value = 2751
$ws.Range("R3").Value = -3215
$ws.Range("S3").Value = 1492
$ws.Range("T3").Value = 2593
$ws.Range("U3").Value = 158
$ws.Range("V3").Value = 11423
$ws.Range("W3").Value = 6.21
$ws.Range("X3").Value = 4.36
$ws.Range("Y3").Value = 4.32
$ws.Range("Z3").Value = 2.37
$ws.Range("AA3").Value = 85.47
$ws.Range("AB3").Value = 29246.24
$ws.Range("AC3").Value = 7340
$ws.Range("AE3").Value = 172128
$ws.Range("AF3").Value = 12.77
$ws.Range("AG3").Value = 920
$ws.Range("AH3").Value = 0.04
$ws.Range("AI3").Value = 12.54
$ws.Range("AJ3").Value = 12372030
$ws.Range("AD3").ClearContents()  # column removed for this row after the fix

# Row 4
$ws.Range("D4").Value = 22642
$ws.Range("E4").Value = 1463
$ws.Range("F4").Value = 1488
$ws.Range("G4").Value = 1260
$ws.Range("H4").Value = 691
$ws.Range("I4").Value = 693
$ws.Range("J4").Value = -2
$ws.Range("K4").Value = 46505
$ws.Range("L4").Value = 22778
$ws.Range("M4").Value = 23727
$ws.Range("N4").Value = 23467
$ws.Range("O4").Value = 260
$ws.Range("P4").Value = 68
$ws.Range("Q4").Value = 1909
$ws.Range("R4").Value = -3413
$ws.Range("S4").Value = 2659
$ws.Range("T4").Value = 3519
$ws.Range("U4").Value = -1609
$ws.Range("V4").Value = 14237
$ws.Range("W4").Value = 6.46
$ws.Range("X4").Value = 3.05
$ws.Range("Y4").Value = 2.96
$ws.Range("Z4").Value = 1.53
$ws.Range("AA4").Value = 96
$ws.Range("AB4").Value = 29827.21
$ws.Range("AC4").Value = 5107
$ws.Range("AE4").Value = 172911
$ws.Range("AF4").Value = 8.359999999999999
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 0.07000000000000001
$ws.Range("AI4").Value = 19.59
$ws.Range("AJ4").Value = 12372030
$ws.Range("AD4").ClearContents()  # column removed for this row after the fix

# Row 5
$ws.Range("D5").Value = 22793
$ws.Range("E5").Value = 754
$ws.Range("F5").Value = 754
$ws.Range("G5").Value = 2293
$ws.Range("H5").Value = 1215
$ws.Range("I5").Value = 1346
$ws.Range("J5").Value = -131
$ws.Range("K5").Value = 34869
$ws.Range("L5").Value = 21741
$ws.Range("M5").Value = 13128
$ws.Range("N5").Value = 13128
$ws.Range("P5").Value = 44
$ws.Range("Q5").Value = 1447
$ws.Range("R5").Value = -1569
$ws.Range("S5").Value = -1531
$ws.Range("T5").Value = 2591
$ws.Range("U5").Value = -1144
$ws.Range("V5").Value = 13565
$ws.Range("W5").Value = 3.31
$ws.Range("X5").Value = 5.33
$ws.Range("Y5").Value = 7.36
$ws.Range("Z5").Value = 2.99
$ws.Range("AA5").Value = 165.61
$ws.Range("AB5").Value = 49030.42
$ws.Range("AC5").Value = 10766
$ws.Range("AD5").Value = 12.34
$ws.Range("AE5").Value = 159898
$ws.Range("AF5").Value = 0.83
$ws.Range("AG5").Value = 3300
$ws.Range("AH5").Value = 0.25
$ws.Range("AI5").Value = 20.13
$ws.Range("AJ5").Value = 7993460
$ws.Range("O5").ClearContents()  # column removed for this row after the fix

# Row 6
$ws.Range("D6").Value = 23463
$ws.Range("E6").Value = 850
$ws.Range("F6").Value = 850
$ws.Range("G6").Value = -689
$ws.Range("H6").Value = -500
$ws.Range("I6").Value = -480
$ws.Range("K6").Value = 33113
$ws.Range("L6").Value = 20739
$ws.Range("M6").Value = 12374
$ws.Range("N6").Value = 12273
$ws.Range("P6").Value = 44
$ws.Range("Q6").Value = 859
$ws.Range("R6").Value = -1410
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = 1387
$ws.Range("U6").Value = -528
$ws.Range("V6").Value = 14137
$ws.Range("W6").Value = 3.62
$ws.Range("X6").Value = -2.13
$ws.Range("Y6").Value = -3.78
$ws.Range("Z6").Value = -1.47
$ws.Range("AA6").Value = 167.6
$ws.Range("AB6").Value = 47577.77
$ws.Range("AC6").Value = -5480
$ws.Range("AD6").Value = -25.55
$ws.Range("AE6").Value = 149481
$ws.Range("AF6").Value = 0.9399999999999999
$ws.Range("AG6").Value = 2700
$ws.Range("AH6").Value = 0.19
$ws.Range("AI6").Value = -46.14
$ws.Range("AJ6").Value = 7993460

# Row 7
$ws.Range("D7").Value = 24718
$ws.Range("E7").Value = 1185
$ws.Range("G7").Value = 479
$ws.Range("H7").Value = -186
$ws.Range("I7").Value = -134
$ws.Range("K7").Value = 33580
$ws.Range("L7").Value = 21202
$ws.Range("M7").Value = 12378
$ws.Range("N7").Value = 12281
$ws.Range("P7").Value = 42
$ws.Range("Q7").Value = 1584
$ws.Range("R7").Value = -1499
$ws.Range("S7").Value = 24
$ws.Range("T7").Value = 1415
$ws.Range("U7").Value = -275
$ws.Range("W7").Value = 4.79
$ws.Range("X7").Value = -0.75
$ws.Range("Y7").Value = -1.1
$ws.Range("Z7").Value = -0.5600000000000001
$ws.Range("AA7").Value = 171.29
$ws.Range("AC7").Value = -1534
$ws.Range("AD7").Value = -80.84
$ws.Range("AE7").Value = 149582
$ws.Range("AF7").Value = 0.83
$ws.Range("AG7").Value = 2438
$ws.Range("AH7").Value = 1.97
$ws.Range("AI7").Value = -144.86

# Row 8
$ws.Range("D8").Value = 25557
$ws.Range("E8").Value = 1319
$ws.Range("G8").Value = 825
$ws.Range("H8").Value = 582
$ws.Range("I8").Value = 618
$ws.Range("K8").Value = 34247
$ws.Range("L8").Value = 21473
$ws.Range("M8").Value = 12776
$ws.Range("N8").Value = 12657
$ws.Range("P8").Value = 42
$ws.Range("Q8").Value = 2241
$ws.Range("R8").Value = -1302
$ws.Range("S8").Value = -514
$ws.Range("T8").Value = 1269
$ws.Range("U8").Value = 941
$ws.Range("W8").Value = 5.16
$ws.Range("X8").Value = 2.28
$ws.Range("Y8").Value = 4.96
$ws.Range("Z8").Value = 1.72
$ws.Range("AA8").Value = 168.06
$ws.Range("AC8").Value = 7049
$ws.Range("AD8").Value = 17.59
$ws.Range("AE8").Value = 154156
$ws.Range("AF8").Value = 0.8
$ws.Range("AG8").Value = 2600
$ws.Range("AH8").Value = 2.1
$ws.Range("AI8").Value = 33.62

# Row 9
$ws.Range("D9").Value = 26530
$ws.Range("E9").Value = 1458
$ws.Range("G9").Value = 1006
$ws.Range("H9").Value = 711
$ws.Range("I9").Value = 683
$ws.Range("K9").Value = 34664
$ws.Range("L9").Value = 21390
$ws.Range("M9").Value = 13274
$ws.Range("N9").Value = 13128
$ws.Range("P9").Value = 42
$ws.Range("Q9").Value = 2199
$ws.Range("R9").Value = -1280
$ws.Range("S9").Value = -564
$ws.Range("T9").Value = 1250
$ws.Range("U9").Value = 891
$ws.Range("W9").Value = 5.5
$ws.Range("X9").Value = 2.68
$ws.Range("Y9").Value = 5.3
$ws.Range("Z9").Value = 2.06
$ws.Range("AA9").Value = 161.14
$ws.Range("AC9").Value = 7792
$ws.Range("AD9").Value = 15.91
$ws.Range("AE9").Value = 159903
$ws.Range("AF9").Value = 0.78
$ws.Range("AG9").Value = 2606
$ws.Range("AH9").Value = 2.1
$ws.Range("AI9").Value = 30.48
